$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CONDUCTOR_operation")

# Insert a new row 10, inheriting formatting from the row above (row 9)
$ws.Rows.Item(10).Insert(-4121, -4163)  # xlShiftDown, xlFormatFromLeftOrAbove

# Populate the new ELECTRIC_SOLVER variable row
$ws.Range("A10").Value = "ELECTRIC_SOLVER"
$ws.Range("B10").Value = "-"
$ws.Range("C10").Value = "integer"
$ws.Range("D10").Value = "Flag to select the solver for the electric module. Possible values: 0= steady state; 1 = transient. Defaults to 1"
$ws.Range("E10").Value = 1

$ws.Rows.Item(10).RowHeight = 101.5

$ws.Application.CutCopyMode = 0

$ws.Select()
$ws.Range("H10").Select()
